# edit.ps1 - applies the "chapter 4" diff:
#   1. Splits the "smartphones" sentence into 3 runs with proofErr spell-check
#      markers around "smartphones".
#   2. Splits the "signale" sentence into 3 runs with proofErr spell-check
#      markers around "signale".
#   3. Appends 4 new bullet (ListParagraph / numId 2) paragraphs at the end
#      of the "Language Learning Models" list.

$d = $word.ActiveDocument

# Unicode characters used below (curly quotes / apostrophe).
$ldq  = [char]0x201C   # “
$rdq  = [char]0x201D   # ”
$apos = [char]0x2019   # '

$pkgHeader = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------------
# 1) "Eventually, the technologies matured ... smartphones ... do."
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Eventually, the technologies matured and gave us everything from smartphones to genetically modified rice. But there were limits to what we could do.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $r.Paragraphs(1)
$pr = $para.Range

$body = '<w:p w:rsidR="00281C4F" w:rsidRDefault="00281C4F" w:rsidP="00281C4F">' +
        '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:jc w:val="both"/></w:pPr>' +
        '<w:r><w:t xml:space="preserve">Eventually, the technologies matured and gave us everything from </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>smartphones</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> to genetically modified rice. But there were limits to what we could do.</w:t></w:r>' +
        '</w:p>'
$xml = $pkgHeader + $body + $pkgFooter
[void]$pr.InsertXML($xml)

# ---------------------------------------------------------------------------
# 2) "The challenge lies in designing an algorithm that ... signale ..."
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("The challenge lies in designing an algorithm that " + $ldq + "knows where to look" + $rdq + " for signale in a given sentence.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $r.Paragraphs(1)
$pr = $para.Range

$body = '<w:p w:rsidR="00583F6C" w:rsidRDefault="00583F6C" w:rsidP="00583F6C">' +
        '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:jc w:val="both"/></w:pPr>' +
        '<w:r><w:t xml:space="preserve">The challenge lies in designing an algorithm that &#8220;knows where to look&#8221; for </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>signale</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> in a given sentence.</w:t></w:r>' +
        '</w:p>'
$xml = $pkgHeader + $body + $pkgFooter
[void]$pr.InsertXML($xml)

# ---------------------------------------------------------------------------
# 3) Append four new bullet paragraphs at the end of the document.
#    The last paragraph in the doc is currently empty (ListParagraph / numId 2)
#    - we fill it with the first new sentence, then clone three more
#    paragraphs after it via InsertParagraphAfter (inherits pPr/numPr/jc).
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$lastPara.Range.Text = "It" + $apos + "s worth noting that humans do this with words of course, but the model doesn" + $apos + "t use our vocabulary."

$lastPara.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.Text = "Instead, it creates a new vocabulary of common tokens that helps it spot patterns across billions of billions of documents"

$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.Text = "In the attention map, every token bears some relationship to every token before it, and for a given input sentence the strength of this relationship describes something about the importance of that token in the sentence."

$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.Text = "In effect, the LLM learns which words to pay attention to."

Write-Output "done"
